# Apply the betexplorer data-correction edit:
# For a handful of rows the "home/away teams + odds + url" block (columns F:V)
# had been written into the wrong physical row. We restore the correct
# F:V block per row by rotating the blocks among the affected rows, and we
# append one new match row (152) that was missing from the sheet.
# Columns A (Indice) and E (data_partida) always stay attached to their own
# physical row, only F:V (the match-specific payload) moves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-FV($row) {
    return $ws.Range("F" + $row + ":V" + $row).Value2
}

function Set-FV($row, $data) {
    $ws.Range("F" + $row + ":V" + $row).Value2 = $data
}

# --- capture all "old" F:V blocks that participate in a rotation ---
$old16  = Get-FV 16
$old18  = Get-FV 18

$old96  = Get-FV 96
$old97  = Get-FV 97

$old102 = Get-FV 102
$old104 = Get-FV 104

$old118 = Get-FV 118
$old120 = Get-FV 120

$old119 = Get-FV 119
$old121 = Get-FV 121

$old122 = Get-FV 122
$old123 = Get-FV 123

$old134 = Get-FV 134
$old135 = Get-FV 135
$old136 = Get-FV 136
$old138 = Get-FV 138

$old146 = Get-FV 146
$old147 = Get-FV 147
$old148 = Get-FV 148
$old149 = Get-FV 149
$old150 = Get-FV 150

# --- simple 2-way swaps ---
Set-FV 16  $old18
Set-FV 18  $old16

Set-FV 96  $old97
Set-FV 97  $old96

Set-FV 102 $old104
Set-FV 104 $old102

Set-FV 118 $old120
Set-FV 120 $old118

Set-FV 119 $old121
Set-FV 121 $old119

Set-FV 122 $old123
Set-FV 123 $old122

# --- 4-way rotation: 134 -> 135 -> 138 -> 136 -> 134 ---
Set-FV 134 $old135
Set-FV 135 $old138
Set-FV 138 $old136
Set-FV 136 $old134

# --- 5-way rotation: 146 -> 147 -> 148 -> 149 -> 150 -> 146 ---
Set-FV 146 $old147
Set-FV 147 $old148
Set-FV 148 $old149
Set-FV 149 $old150
Set-FV 150 $old146

# --- append the new match row (152) that was missing ---
# Copy formatting from the previous last row (151) first, so the new row's
# style ids (bold/border for A, date-time number format for E) line up with
# the ones already used throughout the sheet instead of minting new ones.
$ws.Range("A151:V151").Copy()
$ws.Range("A152:V152").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(152, 1).Value2 = 151
$ws.Cells.Item(152, 2).Value2 = "bulgaria"
$ws.Cells.Item(152, 3).Value2 = "vtora-liga"
$ws.Cells.Item(152, 4).Value2 = "2023-2024"
$ws.Cells.Item(152, 5).Value2 = 45247.5625
$ws.Cells.Item(152, 6).Value2 = "Yantra Gabrovo"
$ws.Cells.Item(152, 7).Value2 = 1
$ws.Cells.Item(152, 8).Value2 = "Ludogorets II"
$ws.Cells.Item(152, 9).Value2 = 1
$ws.Cells.Item(152, 10).Value2 = 1.67
$ws.Cells.Item(152, 11).Value2 = "17/11/2023 02:42"
$ws.Cells.Item(152, 12).Value2 = 1.56
$ws.Cells.Item(152, 13).Value2 = "17/11/2023 13:28"
$ws.Cells.Item(152, 14).Value2 = 3.37
$ws.Cells.Item(152, 15).Value2 = "17/11/2023 02:42"
$ws.Cells.Item(152, 16).Value2 = 3.64
$ws.Cells.Item(152, 17).Value2 = "17/11/2023 13:28"
$ws.Cells.Item(152, 18).Value2 = 4.41
$ws.Cells.Item(152, 19).Value2 = "17/11/2023 02:42"
$ws.Cells.Item(152, 20).Value2 = 5.19
$ws.Cells.Item(152, 21).Value2 = "17/11/2023 13:28"
$ws.Cells.Item(152, 22).Value2 = "https://www.betexplorer.com/football/bulgaria/vtora-liga/yantra-gabrovo-ludogorets/GUDWPqyB/"
